$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INFO")

# Update the description text in A1 (adds the word "footprint")
$ws.Range("A1").Value = "The SubCalc footprint template should be filled out as a flat file with entries in all columns for all rows, then saved as a csv."

# Update the active cell selection on the INFO sheet
$ws.Activate()
$ws.Range("H11").Select()
